$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1179.6
$ws.Range("J38").Value = 70
$ws.Range("L38").Value = 210
$ws.Range("N38").Value = -954
$ws.Range("H129").Value = 777.8333
$ws.Range("I129").Value = 448.54544
$ws.Range("K129").Value = 1345.63632
$ws.Range("M129").Value = 3654.36368
$ws.Range("H132").Value = 2524.4707
$ws.Range("I132").Value = 2524.4707
$ws.Range("K132").Value = 7573.4121
$ws.Range("M132").Value = -5043.4121
$ws.Range("H137").Value = 3796.5535
$ws.Range("I137").Value = 3269.7673
$ws.Range("K137").Value = 9809.3019
$ws.Range("M137").Value = -7259.3019
$ws.Range("H138").Value = 7271.9395
$ws.Range("I138").Value = 7018.9287
$ws.Range("J138").Value = 7340.0576
$ws.Range("K138").Value = 21056.7861
$ws.Range("L138").Value = 22020.1728
$ws.Range("M138").Value = -15916.7861
$ws.Range("N138").Value = -32300.1728
$ws.Range("H141").Value = 2595.3076
$ws.Range("I141").Value = 1914.4
$ws.Range("K141").Value = 5743.200000000001
$ws.Range("M141").Value = -563.2000000000007

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 80000
$ws.Range("J7").Value = 80000
$ws.Range("L7").Value = 80000
$ws.Range("N7").Value = -80228
$ws.Range("H32").Value = 35157.957
$ws.Range("I32").Value = 26515.422
$ws.Range("K32").Value = 26515.422
$ws.Range("M32").Value = -26228.422
$ws.Range("H45").Value = 3844.6
$ws.Range("I45").Value = 3078.3333
$ws.Range("J45").Value = 4994
$ws.Range("K45").Value = 3078.3333
$ws.Range("L45").Value = 4994
$ws.Range("M45").Value = -2701.3333
$ws.Range("N45").Value = -5748
$ws.Range("H74").Value = 4778.357
$ws.Range("I74").Value = 3971
$ws.Range("J74").Value = 7200.4287
$ws.Range("K74").Value = 3971
$ws.Range("L74").Value = 7200.4287
$ws.Range("M74").Value = -3097
$ws.Range("N74").Value = -8948.4287
$ws.Range("H77").Value = 4778.357
$ws.Range("I77").Value = 3971
$ws.Range("J77").Value = 7200.4287
$ws.Range("K77").Value = 19855
$ws.Range("L77").Value = 36002.14350000001
$ws.Range("M77").Value = -15487
$ws.Range("N77").Value = -44738.14350000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13898.9375
$ws.Range("I86").Value = 11282.083
$ws.Range("K86").Value = 11282.083
$ws.Range("M86").Value = -10159.083
$ws.Range("H89").Value = 13898.9375
$ws.Range("I89").Value = 11282.083
$ws.Range("K89").Value = 56410.415
$ws.Range("M89").Value = -50794.415
$ws.Range("H107").Value = 3237.125
$ws.Range("I107").Value = 3313.8572
$ws.Range("K107").Value = 3313.8572
$ws.Range("M107").Value = -1393.8572

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7714.696
$ws.Range("I31").Value = 4385.5835
$ws.Range("J31").Value = 19699.5
$ws.Range("K31").Value = 4385.5835
$ws.Range("L31").Value = 19699.5
$ws.Range("M31").Value = -4090.5835
$ws.Range("N31").Value = -20289.5
$ws.Range("H33").Value = 6644.2
$ws.Range("I33").Value = 6644.2
$ws.Range("K33").Value = 6644.2
$ws.Range("M33").Value = -6265.2
$ws.Range("H34").Value = 7714.696
$ws.Range("I34").Value = 4385.5835
$ws.Range("J34").Value = 19699.5
$ws.Range("K34").Value = 4385.5835
$ws.Range("L34").Value = 19699.5
$ws.Range("M34").Value = -4183.5835
$ws.Range("N34").Value = -20103.5
$ws.Range("H99").Value = 5333.3228
$ws.Range("I99").Value = 6494.9546
$ws.Range("J99").Value = 2493.7778
$ws.Range("K99").Value = 6494.9546
$ws.Range("L99").Value = 2493.7778
$ws.Range("M99").Value = -4996.9546
$ws.Range("N99").Value = -5489.7778
$ws.Range("H126").Value = 5333.3228
$ws.Range("I126").Value = 6494.9546
$ws.Range("J126").Value = 2493.7778
$ws.Range("K126").Value = 19484.8638
$ws.Range("L126").Value = 7481.3334
$ws.Range("M126").Value = -17014.8638
$ws.Range("N126").Value = -12421.3334
$ws.Range("H132").Value = 3587.8235
$ws.Range("J132").Value = 5828.5713
$ws.Range("L132").Value = 17485.7139
$ws.Range("N132").Value = -22545.7139
$ws.Range("H134").Value = 3849.5676
$ws.Range("I134").Value = 3511.3076
$ws.Range("J134").Value = 4649.091
$ws.Range("K134").Value = 10533.9228
$ws.Range("L134").Value = 13947.273
$ws.Range("M134").Value = -7998.9228
$ws.Range("N134").Value = -19017.273
$ws.Range("H141").Value = 284248.53
$ws.Range("J141").Value = 335945.62
$ws.Range("L141").Value = 335945.62
$ws.Range("N141").Value = -346305.62

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 7503
$ws.Range("I97").Value = 8755
$ws.Range("J97").Value = 4999
$ws.Range("K97").Value = 8755
$ws.Range("L97").Value = 4999
$ws.Range("M97").Value = -8259
$ws.Range("N97").Value = -5991
$ws.Range("H102").Value = 3381.5264
$ws.Range("I102").Value = 3026.724
$ws.Range("K102").Value = 3026.724
$ws.Range("M102").Value = -1404.724
$ws.Range("H107").Value = 439.58334
$ws.Range("I107").Value = 459.6
$ws.Range("K107").Value = 459.6
$ws.Range("M107").Value = 1460.4
$ws.Range("H113").Value = 1619.375
$ws.Range("I113").Value = 1635.5
$ws.Range("J113").Value = 1571
$ws.Range("K113").Value = 1635.5
$ws.Range("L113").Value = 1571
$ws.Range("M113").Value = 534.5
$ws.Range("N113").Value = -5911
$ws.Range("H122").Value = 4391.615
$ws.Range("I122").Value = 2781
$ws.Range("K122").Value = 8343
$ws.Range("M122").Value = -5893
$ws.Range("H126").Value = 10468.806
$ws.Range("I126").Value = 9926.931
$ws.Range("J126").Value = 12713.714
$ws.Range("K126").Value = 29780.793
$ws.Range("L126").Value = 38141.142
$ws.Range("M126").Value = -27310.793
$ws.Range("N126").Value = -43081.142
$ws.Range("H132").Value = 7669.3887
$ws.Range("I132").Value = 6959.6514
$ws.Range("K132").Value = 20878.9542
$ws.Range("M132").Value = -18348.9542

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2494.5
$ws.Range("I7").Value = 989
$ws.Range("K7").Value = 989
$ws.Range("M7").Value = -877
$ws.Range("H40").Value = 5734.5
$ws.Range("I40").Value = 5082.273
$ws.Range("K40").Value = 5082.273
$ws.Range("M40").Value = -4946.273
$ws.Range("H46").Value = 3914.4707
$ws.Range("J46").Value = 3789.3
$ws.Range("L46").Value = 3789.3
$ws.Range("N46").Value = -4165.3
$ws.Range("H82").Value = 1112.375
$ws.Range("I82").Value = 675
$ws.Range("J82").Value = 1549.75
$ws.Range("K82").Value = 675
$ws.Range("L82").Value = 1549.75
$ws.Range("M82").Value = -314
$ws.Range("N82").Value = -2271.75
$ws.Range("H85").Value = 1112.375
$ws.Range("I85").Value = 675
$ws.Range("J85").Value = 1549.75
$ws.Range("K85").Value = 675
$ws.Range("L85").Value = 1549.75
$ws.Range("M85").Value = 573
$ws.Range("N85").Value = -4045.75
$ws.Range("H122").Value = 5544.722
$ws.Range("I122").Value = 4343.2144
$ws.Range("K122").Value = 13029.6432
$ws.Range("M122").Value = -10579.6432
$ws.Range("H126").Value = 2494.5
$ws.Range("I126").Value = 989
$ws.Range("K126").Value = 2967
$ws.Range("M126").Value = -497
$ws.Range("H132").Value = 4325.3335
$ws.Range("I132").Value = 3772.2593
$ws.Range("J132").Value = 6814.1665
$ws.Range("K132").Value = 11316.7779
$ws.Range("L132").Value = 20442.4995
$ws.Range("M132").Value = -8786.777900000001
$ws.Range("N132").Value = -25502.4995

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2384.8857
$ws.Range("I126").Value = 1876.9656
$ws.Range("K126").Value = 5630.8968
$ws.Range("M126").Value = -3160.8968
$ws.Range("H132").Value = 6674.9443
$ws.Range("I132").Value = 6720
$ws.Range("J132").Value = 6449.6665
$ws.Range("K132").Value = 20160
$ws.Range("L132").Value = 19348.9995
$ws.Range("M132").Value = -17630
$ws.Range("N132").Value = -24408.9995
$ws.Range("H136").Value = 5370.852
$ws.Range("I136").Value = 4793.05
$ws.Range("K136").Value = 14379.15
$ws.Range("M136").Value = -11829.15
